$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "confirmed"
$ws.Range("H3").Value = "confirmed"
$ws.Range("H7").Value = "cancelled"
